# Fruta / hortaliza, semanal
# Insert a new data row at row 409 (pushing existing rows 409.. down by one)
# and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 409; this shifts rows 409..488 -> 410..489
$ws.Rows.Item(409).Insert()

# Copy the number format of the date column from the row above so the new
# row's date cell matches the existing date formatting (style s="2").
$ws.Cells.Item(408, 4).Copy()
$ws.Cells.Item(409, 4).PasteSpecial(-4122) # xlPasteFormats

# Populate the new row's values
$ws.Cells.Item(409, 1).Value = 5
$ws.Cells.Item(409, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(409, 3).Value = "Maule"
$ws.Cells.Item(409, 4).Value = 44995
$ws.Cells.Item(409, 5).Value = 7
$ws.Cells.Item(409, 6).Value = 100114013
$ws.Cells.Item(409, 7).Value = "Zanahoria"
$ws.Cells.Item(409, 8).Value = "Sin especificar"
$ws.Cells.Item(409, 9).Value = "Primera"
$ws.Cells.Item(409, 10).Value = 400
$ws.Cells.Item(409, 11).Value = 7000
$ws.Cells.Item(409, 12).Value = 7000
$ws.Cells.Item(409, 13).Value = 7000
$ws.Cells.Item(409, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(409, 15).Value = "Región de Ñuble"
$ws.Cells.Item(409, 16).Value = 350
$ws.Cells.Item(409, 17).Value = 20
$ws.Cells.Item(409, 18).Value = "Hortaliza"
